$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows that were entirely removed (descending order to keep row numbers stable)
$ws.Rows(28).EntireRow.Delete()
$ws.Rows(26).EntireRow.Delete()

# Apply individual cell value changes (row numbers refer to post-deletion layout)
$ws.Range("C3").Value = 11.2
$ws.Range("E4").ClearContents()
$ws.Range("C5").ClearContents()
$ws.Range("F6").Value = 16.43
$ws.Range("E9").Value = -6.8
$ws.Range("E10").Value = -6.1
$ws.Range("F12").ClearContents()
$ws.Range("F14").Value = 17.76
$ws.Range("E17").ClearContents()
$ws.Range("F17").Value = 17.78
$ws.Range("E18").ClearContents()
$ws.Range("F19").Value = 17.81
$ws.Range("F20").ClearContents()
$ws.Range("C21").Value = 12.7
$ws.Range("C23").ClearContents()
$ws.Range("F23").ClearContents()
$ws.Range("F25").ClearContents()
$ws.Range("F27").Value = 17
$ws.Range("F28").Value = 17.44
$ws.Range("C32").Value = 10.5
$ws.Range("D32").Value = -14.7
